# Adds two new weekly price records to the "Repollo" sheet.
# 1) A new row is inserted at row 365 (pushing the former rows 365-461 down by one).
# 2) A second new row is inserted at row 460 -- i.e. right before the former row 459
#    (now sitting at row 460 after the first insertion) -- pushing it (and the rows
#    after it) down by one more, so the sheet ends with 462 data rows (+ header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First new record --------------------------------------------------
$ws.Rows.Item(365).EntireRow.Insert()

$ws.Range("A365").Value = 7
$ws.Range("B365").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C365").Value = "Ñuble"
$ws.Range("D365").Value = 45120
$ws.Range("E365").Value = 16
$ws.Range("F365").Value = 100112006
$ws.Range("G365").Value = "Repollo"
$ws.Range("H365").Value = "Crespo record"
$ws.Range("I365").Value = "Primera"
$ws.Range("J365").Value = 250
$ws.Range("K365").Value = 1000
$ws.Range("L365").Value = 1000
$ws.Range("M365").Value = 1000
$ws.Range("N365").Value = "`$/unidad"
$ws.Range("O365").Value = "Provincia de Diguillín"
$ws.Range("P365").Value = 1000
$ws.Range("Q365").Value = 1
$ws.Range("R365").Value = "Hortaliza"

# --- Second new record ---------------------------------------------------
$ws.Rows.Item(460).EntireRow.Insert()

$ws.Range("A460").Value = 7
$ws.Range("B460").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C460").Value = "Ñuble"
$ws.Range("D460").Value = 45121
$ws.Range("E460").Value = 16
$ws.Range("F460").Value = 100112006
$ws.Range("G460").Value = "Repollo"
$ws.Range("H460").Value = "Crespo record"
$ws.Range("I460").Value = "Primera"
$ws.Range("J460").Value = 250
$ws.Range("K460").Value = 1000
$ws.Range("L460").Value = 1000
$ws.Range("M460").Value = 1000
$ws.Range("N460").Value = "`$/unidad"
$ws.Range("O460").Value = "Provincia de Diguillín"
$ws.Range("P460").Value = 1000
$ws.Range("Q460").Value = 1
$ws.Range("R460").Value = "Hortaliza"
